$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.961.72'
$ws.Range("E2").Value = '  -1.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.468.79'
$ws.Range("E3").Value = '  -2.89%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.21'
$ws.Range("E5").Value = '  -1.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.48'
$ws.Range("E6").Value = '  -3.25%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.514'
$ws.Range("E8").Value = '  -2.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.468.64'
$ws.Range("E9").Value = '  -2.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.134'
$ws.Range("E10").Value = '  -3.16%  '
$ws.Range("E11").Value = '  -1.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.92'
$ws.Range("E12").Value = '  -2.76%  '
$ws.Range("E13").Value = '  -3.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.54'
$ws.Range("E14").Value = '  -3.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.924.87'
$ws.Range("E15").Value = '  -1.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.591.45'
$ws.Range("E16").Value = '  -2.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000169'
$ws.Range("E17").Value = '  -4.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.418.62'
$ws.Range("E18").Value = '  -4.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.07'
$ws.Range("E19").Value = '  -6.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.52'
$ws.Range("E20").Value = '  -6.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '351.76'
$ws.Range("E21").Value = '  -5.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.03'
$ws.Range("E22").Value = '  -2.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.80'
$ws.Range("E24").Value = '  -4.42%  '
$ws.Range("B25").Value = 'NEARProtocol'
$ws.Range("C25").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.24'
$ws.Range("E25").Value = '  -7.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.82'
$ws.Range("E26").Value = '  -5.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.18'
$ws.Range("E27").Value = '  -7.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  -57.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.596.38'
$ws.Range("E29").Value = '  -2.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0903'
$ws.Range("E30").Value = '  -7.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '508.67'
$ws.Range("E31").Value = '  -5.52%  '
$ws.Range("E32").Value = '  -7.59%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.77'
$ws.Range("E33").Value = '  -5.50%  '
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.23'
$ws.Range("E34").Value = '  -7.12%  '
$ws.Range("E35").Value = '  +0.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '159.61'
$ws.Range("E36").Value = '  +0.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.115'
$ws.Range("E37").Value = '  -10.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.65'
$ws.Range("E38").Value = '  +0.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.33'
$ws.Range("E39").Value = '  -4.72%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.33'
$ws.Range("E40").Value = '  -7.89%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.69'
$ws.Range("E41").Value = '  -5.47%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.82'
$ws.Range("E43").Value = '  -6.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.327'
$ws.Range("E44").Value = '  -6.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.37'
$ws.Range("E45").Value = '  -6.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.95'
$ws.Range("E46").Value = '  -1.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '140.76'
$ws.Range("E47").Value = '  -4.99%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.514'
$ws.Range("E48").Value = '  -6.82%  '
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.44'
$ws.Range("E49").Value = '  -7.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0255'
$ws.Range("E50").Value = '  -10.90%  '
$ws.Range("E51").Value = '  -7.49%  '
